$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

# Stash the cell's current formatting (borders/fill/font/number format) in an
# unused scratch cell, since assigning a new value further below can cause
# the engine to mint a fresh style record (e.g. when the literal text looks
# like a number, Excel marks the cell with a quote-prefix / text number
# format). Copying the format back afterwards restores the original look.
$scratch = $ws.Range("Z1")
$cell.Copy($scratch)

# B11 held the text "R40" (row label); it should now read "1".
$cell.Value = "'1"

# Restore the original formatting only (leave the new value untouched).
$scratch.Copy()
$cell.PasteSpecial(-4122)

# Clean up the scratch cell.
$scratch.Clear()
